$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.688.48"
$ws.Range("E2").Value = "  -6.58%  "

$ws.Range("D3").Value = "3.529.20"
$ws.Range("E3").Value = "  -2.67%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "395.56"
$ws.Range("E5").Value = "  -5.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "124.96"
$ws.Range("E6").Value = "  -4.96%  "

$ws.Range("D7").Value = "3.523.26"
$ws.Range("E7").Value = "  -2.68%  "

$ws.Range("E8").Value = "  -8.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.685"
$ws.Range("E10").Value = "  -11.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.153"
$ws.Range("E11").Value = "  -16.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000340"
$ws.Range("E12").Value = "  -4.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.29"
$ws.Range("E13").Value = "  -8.01%  "

$ws.Range("D14").Value = "4.085.06"
$ws.Range("E14").Value = "  -2.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.29"
$ws.Range("E15").Value = "  -6.52%  "

$ws.Range("E16").Value = "  -3.12%  "

$ws.Range("D17").Value = "3.514.52"
$ws.Range("E17").Value = "  -3.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.87"
$ws.Range("E18").Value = "  -7.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.69"
$ws.Range("E19").Value = "  +1.96%  "

$ws.Range("D20").Value = "63.777.44"
$ws.Range("E20").Value = "  -6.36%  "

$ws.Range("E21").Value = "  -9.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "398.50"
$ws.Range("E22").Value = "  -14.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.99"
$ws.Range("E23").Value = "  +4.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.61"
$ws.Range("E24").Value = "  -8.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.89"
$ws.Range("E25").Value = "  -7.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "34.02"
$ws.Range("E26").Value = "  -5.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.23"
$ws.Range("E27").Value = "  +7.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.01"
$ws.Range("E28").Value = "  -10.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.85"
$ws.Range("E29").Value = "  -12.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.95"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.57"
$ws.Range("E31").Value = "  -7.33%  "

$ws.Range("E32").Value = "  -5.31%  "

$ws.Range("E33").Value = "  -7.00%  "

$ws.Range("E34").Value = "  -7.06%  "

$ws.Range("E35").Value = "  +0.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "36.96"
$ws.Range("E36").Value = "  -9.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.35"
$ws.Range("E37").Value = "  -4.33%  "

$ws.Range("E38").Value = "  -10.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("E40").Value = "  +19.83%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.132"
$ws.Range("E41").Value = "  -9.42%  "

$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0629"
$ws.Range("E42").Value = "  -11.33%  "

$ws.Range("E43").Value = "  +13.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "140.72"
$ws.Range("E44").Value = "  -5.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.10"
$ws.Range("E45").Value = "  -5.25%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.97"
$ws.Range("E46").Value = "  +14.92%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.73"
$ws.Range("E47").Value = "  -10.36%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.10"
$ws.Range("E48").Value = "  -5.29%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.95"
$ws.Range("E49").Value = "  -1.39%  "

$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.48"
$ws.Range("E50").Value = "  -9.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.277"
$ws.Range("E51").Value = "  -10.26%  "
